$d = $word.ActiveDocument

# --- Paragraph 1 (Heading2 "Supplement Table S2"): ---
#   run 1: "Supplement"   -> "SUPPLEMENT TABLE S"
#   run 2: " Table S2"    -> "2"
# and the Word-managed "_GoBack" bookmark (which marks the location of the
# last text edit) moves from the end of the document to the boundary
# between these two runs.

$oldRun1 = "Supplement"
$newRun1 = "SUPPLEMENT TABLE S"
$oldRun2 = " Table S2"
$newRun2 = "2"

$p1 = $d.Paragraphs(1)
$paraStart = $p1.Range.Start

# Insert the _GoBack bookmark at the boundary between run 1 and run 2 BEFORE
# editing any text. Doing this first keeps the two runs from being merged
# back into one despite matching formatting (the bookmark sits between
# them), and it also removes the old "_GoBack" bookmark wherever it
# currently is in the document, since Word keeps only one instance of it
# at a time (it always marks the most recent edit position).
$boundary = $d.Range($paraStart + $oldRun1.Length, $paraStart + $oldRun1.Length)
$d.Bookmarks.Add("_GoBack", $boundary)

# Edit run 1's text in place.
$run1Range = $d.Range($paraStart, $paraStart + $oldRun1.Length)
$run1Range.Text = $newRun1

# Run 2 now begins right after the replacement text for run 1.
$run2Start = $paraStart + $newRun1.Length
$run2Range = $d.Range($run2Start, $run2Start + $oldRun2.Length)
$run2Range.Text = $newRun2

Write-Output ("Paragraph 1 now reads: " + $d.Paragraphs(1).Range.Text)
